$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.051938
$ws.Range("H2").Value = 0.155814
$ws.Range("I2").Value = 0.1172837182974765
$ws.Range("J2").Value = 0.1172837182974765
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 19.523753
$ws.Range("N2").Value = 58.571259
$ws.Range("O2").Value = 0.4652898160202426
$ws.Range("P2").Value = 0.4652898160202426
$ws.Range("Q2").Value = 1.014024683314
$ws.Range("R2").Value = 9.126222149826001
$ws.Range("S2").Value = 0.05457091970880278
$ws.Range("T2").Value = 0.05457091970880278

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.051938
$ws.Range("H3").Value = 0.155814
$ws.Range("I3").Value = 0.1172837182974765
$ws.Range("J3").Value = 0.1172837182974765
$ws.Range("O3").Value = 0.4569298967820781
$ws.Range("P3").Value = 0.4569298967820781
$ws.Range("Q3").Value = 0.9958055773586667
$ws.Range("R3").Value = 8.962250196228
$ws.Range("S3").Value = 0.05359043729588424
$ws.Range("T3").Value = 0.05359043729588425

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.051938
$ws.Range("H4").Value = 0.155814
$ws.Range("I4").Value = 0.1172837182974765
$ws.Range("J4").Value = 0.1172837182974765
$ws.Range("O4").Value = 0.07778028719767933
$ws.Range("P4").Value = 0.07778028719767933
$ws.Range("Q4").Value = 0.169509687034
$ws.Range("R4").Value = 1.525587183306
$ws.Range("S4").Value = 0.009122361292789437
$ws.Range("T4").Value = 0.009122361292789437

$ws.Range("I5").Value = 0.5474768201053503
$ws.Range("J5").Value = 0.5474768201053503
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 19.523753
$ws.Range("N5").Value = 58.571259
$ws.Range("O5").Value = 0.4652898160202426
$ws.Range("P5").Value = 0.4652898160202426
$ws.Range("Q5").Value = 4.733436296085
$ws.Range("R5").Value = 42.60092666476501
$ws.Range("S5").Value = 0.2547353889021658
$ws.Range("T5").Value = 0.2547353889021659

$ws.Range("I6").Value = 0.5474768201053503
$ws.Range("J6").Value = 0.5474768201053503
$ws.Range("O6").Value = 0.4569298967820781
$ws.Range("P6").Value = 0.4569298967820781
$ws.Range("S6").Value = 0.250158526901318
$ws.Range("T6").Value = 0.2501585269013181

$ws.Range("I7").Value = 0.5474768201053503
$ws.Range("J7").Value = 0.5474768201053503
$ws.Range("O7").Value = 0.07778028719767933
$ws.Range("P7").Value = 0.07778028719767933
$ws.Range("S7").Value = 0.04258290430186637
$ws.Range("T7").Value = 0.04258290430186637

$ws.Range("I8").Value = 0.3352394615971734
$ws.Range("J8").Value = 0.3352394615971734
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 19.523753
$ws.Range("N8").Value = 58.571259
$ws.Range("O8").Value = 0.4652898160202426
$ws.Range("P8").Value = 0.4652898160202426
$ws.Range("Q8").Value = 2.898450814956333
$ws.Range("R8").Value = 26.086057334607
$ws.Range("S8").Value = 0.155983507409274
$ws.Range("T8").Value = 0.155983507409274

$ws.Range("I9").Value = 0.3352394615971734
$ws.Range("J9").Value = 0.3352394615971734
$ws.Range("O9").Value = 0.4569298967820781
$ws.Range("P9").Value = 0.4569298967820781
$ws.Range("S9").Value = 0.1531809325848759
$ws.Range("T9").Value = 0.1531809325848759

$ws.Range("I10").Value = 0.3352394615971734
$ws.Range("J10").Value = 0.3352394615971734
$ws.Range("O10").Value = 0.07778028719767933
$ws.Range("P10").Value = 0.07778028719767933
$ws.Range("Q10").Value = 0.4845202474963334
$ws.Range("R10").Value = 4.360682227467001
$ws.Range("S10").Value = 0.02607502160302354
$ws.Range("T10").Value = 0.02607502160302354
